# Re-shuffle the per-trial detail columns (category, cond_cat, correct_answer,
# stimulus, conceptual, perceptual, typicality, n, p_*, r_*) across data rows
# 2-41, while leaving the trial-index columns (A-G, J) untouched.
#
# $rowMap[destRow] = srcRow  -- destRow's detail block becomes a copy of the
# *original* (pre-edit) srcRow's detail block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2 = 8
    3 = 23
    4 = 38
    5 = 31
    6 = 14
    7 = 18
    8 = 26
    9 = 34
    10 = 13
    11 = 39
    12 = 37
    13 = 35
    14 = 30
    15 = 27
    16 = 17
    17 = 24
    18 = 2
    19 = 3
    20 = 19
    21 = 41
    22 = 29
    23 = 16
    24 = 40
    25 = 33
    26 = 11
    27 = 22
    28 = 6
    29 = 7
    30 = 32
    31 = 15
    32 = 10
    33 = 12
    34 = 9
    35 = 5
    36 = 36
    37 = 28
    38 = 20
    39 = 21
    40 = 25
    41 = 4
}

# Columns whose values travel together as one block per source row.
# (1-indexed spreadsheet columns: H=8 I=9 K=11 L=12 M=13 N=14 O=15 P=16
#  Q=17 R=18 S=19 T=20 U=21 V=22)
$blockCols = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)

# Snapshot every source row's block BEFORE any writes happen, since some
# destination rows read from rows that are themselves about to be overwritten.
# NOTE: use .Value2 (not .Value) -- in this host, round-tripping a bare
# .Value read through a PowerShell variable loses the underlying scalar, so
# every read/write below goes through .Value2 instead.
$snapshot = @{}
foreach ($r in 2..41) {
    $rowVals = @{}
    foreach ($c in $blockCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in 2..41) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $blockCols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
